$wb = $excel.ActiveWorkbook

# --- 1. Add the two new summary formulas on the "Lat Long by Individual"
#        sheet (combined average for populations 1 and 2, row 109) ---
$wsIndiv = $wb.Worksheets.Item("Lat Long by Individual")
$wsIndiv.Range("E109").Formula = "=AVERAGE(C106:C125, C127:C147)"
$wsIndiv.Range("F109").Formula = "=AVERAGE(D106:D125, D127:D147)"

# --- 2. Add a new worksheet "PRBI_Comb12_AvgLatLong" after the existing
#        "PRBI_AvgLatLong" sheet, with the combined-population averages ---
$wsLast = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsLast)
$ws3.Name = "PRBI_Comb12_AvgLatLong"

$ws3.Range("A1").Value = "Population"
$ws3.Range("B1").Value = "Avg. Latitude"
$ws3.Range("C1").Value = "Avg. Longitude"

$ws3.Range("A2").Value = "1_2"
$ws3.Range("B2").Value = 10.201873197560973
$ws3.Range("C2").Value = 124.21258387560971

$data = @(
    @(7, 9.414770978, 123.31799460000001),
    @(8, 9.616869822, 123.47946159999999),
    @(9, 9.848377696, 123.57853249999999),
    @(10, 10.073306219999999, 123.65455710000001),
    @(11, 10.224988701111108, 123.80282073333332),
    @(19, 10.01649883, 125.02510150000001)
)

$r = 3
foreach ($row in $data) {
    $ws3.Cells.Item($r, 1).Value = $row[0]
    $ws3.Cells.Item($r, 2).Value = $row[1]
    $ws3.Cells.Item($r, 3).Value = $row[2]
    $r++
}

# Make the new sheet the active/selected tab, matching the authored workbook.
$ws3.Activate()
